$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.857.80'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.917.49'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.32'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3807'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07763'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9767'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.940.76'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.982'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.698'
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06976'
$ws.Range('E15').Value = '  -1.17%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.005'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '84.41'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000009497'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.64'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '28.846.05'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.342'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.10'
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.149.49'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.065'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '157.89'
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.98'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.621'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '118.00'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.847'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09316'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.8722'
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.101'
$ws.Range('E33').Value = '  +0.67%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.247'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.031'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.150'
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02038'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('E40').Value = '  +11.34%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.512'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5500'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1756'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.000002888'
$ws.Range('E44').Value = '  +16.79%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '9.341'
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.174'
$ws.Range('E46').Value = '  +3.96%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5164'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('E49').Value = '  -1.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '110.62'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.764'
$ws.Range('E51').Value = '  -0.65%  '
